$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.692.57'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '3.077.46'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("D4").Value = '''1.01'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '''234.37'
$ws.Range("E5").Value = '  +8.90%  '
$ws.Range("D6").Value = '''616.23'
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("E7").Value = '  -10.72%  '
$ws.Range("D8").Value = '''0.359'
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").Value = '3.074.75'
$ws.Range("E10").Value = '  -2.30%  '
$ws.Range("D11").Value = '''0.710'
$ws.Range("E11").Value = '  -5.67%  '
$ws.Range("D12").Value = '''0.197'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").Value = '''0.0000249'
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("D14").Value = '''35.07'
$ws.Range("E14").Value = '  +0.75%  '
$ws.Range("D15").Value = '89.717.16'
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").Value = '''5.35'
$ws.Range("E16").Value = '  -6.56%  '
$ws.Range("D17").Value = '3.669.83'
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.154.73'
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("B19").Value = 'SuiNetwork'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").Value = '''3.80'
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").Value = '''0.0000209'
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = '''13.75'
$ws.Range("E21").Value = '  -5.77%  '
$ws.Range("D22").Value = '''431.96'
$ws.Range("E22").Value = '  -9.08%  '
$ws.Range("D23").Value = '''5.39'
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("D24").Value = '''8.74'
$ws.Range("E24").Value = '  -4.12%  '
$ws.Range("D25").Value = '''5.56'
$ws.Range("E25").Value = '  -2.99%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '''11.70'
$ws.Range("E26").Value = '  -4.81%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '''81.32'
$ws.Range("E27").Value = '  -15.27%  '
$ws.Range("D28").Value = '3.285.81'
$ws.Range("E28").Value = '  -1.39%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''9.01'
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").Value = '''0.157'
$ws.Range("E32").Value = '  -3.35%  '
$ws.Range("D33").Value = '''0.191'
$ws.Range("E33").Value = '  -5.89%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '''0.152'
$ws.Range("E34").Value = '  +5.51%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '''25.53'
$ws.Range("E35").Value = '  -6.96%  '
$ws.Range("D36").Value = '''3.67'
$ws.Range("E36").Value = '  +2.33%  '
$ws.Range("D37").Value = '''7.04'
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("D38").Value = '''492.56'
$ws.Range("E38").Value = '  -4.78%  '
$ws.Range("D39").Value = '''1.87'
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("D40").Value = '''1.25'
$ws.Range("E40").Value = '  -3.55%  '
$ws.Range("D41").Value = '''3.58'
$ws.Range("E41").Value = '  +56.57%  '
$ws.Range("D42").Value = '''0.0872'
$ws.Range("E42").Value = '  -4.44%  '
$ws.Range("D43").Value = '''22.12'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '''0.394'
$ws.Range("E45").Value = '  -5.96%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").Value = '''151.88'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '''1.84'
$ws.Range("E47").Value = '  -6.56%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''0.672'
$ws.Range("E48").Value = '  -6.39%  '
$ws.Range("D49").Value = '''44.35'
$ws.Range("E49").Value = '  -2.43%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").Value = '''1.29'
$ws.Range("E51").Value = '  -5.01%  '
